$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 11 (shifts existing rows 11-64 down to 12-65),
# copying format from the row above (row 10), matching Excel's default
# Insert behaviour.
$ws.Rows("11:11").Insert()

# Make sure the freshly inserted row keeps the same (default) row height
# as its neighbours in the "image" block.
$ws.Rows(11).RowHeight = 15.75

# Populate the new row with the Bg_City.jpg asset entry.
$ws.Range("A11").Value = "image"
$ws.Range("B11").Value = "Bg_City.jpg"
$ws.Range("C11").Value = "City BG for panning."
$ws.Range("D11").Value = "3840 X 1080"
$ws.Range("E11").Value = "Image by Kyle"
$ws.Range("F11").Value = "Placeholder"

# F11 should use the same style as the other "Item_*" rows below it
# (style index 2) rather than the style inherited from the row above.
$ws.Range("F11").Style = $ws.Range("F19").Style

# Restore the selection to what the author left it at.
$ws.Range("C6").Select()
